# Progresso trabalho.xlsx -- apply "Add files via upload" edit
#
# Summary of the real content change (the cellXfs/style-index churn visible
# in the raw XML diff is just Excel/LO renumbering style records on save --
# every *effective* format for the pre-existing rows 2-27 is unchanged):
#   1) Column A "done" flags for rows 8-17,19,22,23 go from 1 -> 0
#      (re-opening those sub-tasks as not-yet-done).
#   2) K6's SUM widens to cover the newly reopened rows + the new section.
#   3) K4's percentage-of-completion formula now divides by 25 (was 20).
#   4) A brand-new "5) Relatorio" section is appended at rows 29-34 with
#      5 sub-items (Capa / Indice / Introducao / Implementacao / Conclusao).
#   5) The active-cell selection marker moves to I15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Reset the "done" checkbox column (A) for the reopened sub-tasks.
# ---------------------------------------------------------------------
$reopened = @(8,9,10,11,12,13,14,15,16,17,19,22,23)
foreach ($r in $reopened) {
    $ws.Cells.Item($r, 1).Value = 0
}

# ---------------------------------------------------------------------
# 2) Update the two summary formulas.
# ---------------------------------------------------------------------
$ws.Range("K6").Formula = "=SUM(A4:A5,A8:A24,A27,A30:A34)"
$ws.Range("K4").Formula = "=K6/25"

# ---------------------------------------------------------------------
# 3) Append the new "5) Relatorio" block (rows 29-34).
# ---------------------------------------------------------------------

# Row 29 - section header "5) | Relatorio" (plain boxed cells, no fill)
$ws.Range("B29").Value = "5)"
$ws.Range("C29").Value = "Relatório"
$ws.Range("B29:C29").Borders.LineStyle = 1

# Row 30 - "5).1) | Capa" ; A30 flag, B30 open-top box, C30 full box,
# D30 plain left-aligned cell (no border, not merged with C30)
$ws.Cells.Item(30, 1).Value = 0
$ws.Range("B30").Value = "5).1)"
$ws.Range("C30").Value = "Capa"
$ws.Range("D30").Value = $null

$ws.Range("B30").Borders.Item(7).LineStyle = 1   # xlEdgeLeft
$ws.Range("B30").Borders.Item(10).LineStyle = 1  # xlEdgeRight
$ws.Range("B30").Borders.Item(9).LineStyle = 1   # xlEdgeBottom
$ws.Range("B30").Borders.Item(8).LineStyle = 0   # xlEdgeTop (none)

$ws.Range("C30").Borders.LineStyle = 1
$ws.Range("C30").HorizontalAlignment = -4131     # xlLeft

$ws.Range("D30").HorizontalAlignment = -4131     # xlLeft

# Row 31 - "5).2) | Indice" ; B31 full box, C31 box without bottom edge,
# D31 plain left-aligned cell (no border, not merged with C31)
$ws.Cells.Item(31, 1).Value = 0
$ws.Range("B31").Value = "5).2)"
$ws.Range("C31").Value = "Índice"
$ws.Range("D31").Value = $null

$ws.Range("B31").Borders.LineStyle = 1

$ws.Range("C31").Borders.Item(7).LineStyle = 1   # xlEdgeLeft
$ws.Range("C31").Borders.Item(10).LineStyle = 1  # xlEdgeRight
$ws.Range("C31").Borders.Item(8).LineStyle = 1   # xlEdgeTop
$ws.Range("C31").Borders.Item(9).LineStyle = 0   # xlEdgeBottom (none)
$ws.Range("C31").HorizontalAlignment = -4131     # xlLeft

$ws.Range("D31").HorizontalAlignment = -4131     # xlLeft

# Rows 32-34 - "5).3) Introducao", "5).4) Implementacao", "5).5) Conclusao"
# B full box; C:D merged, full box, left aligned.
$rows32to34 = @(
    @{ Row = 32; Num = "5).3)"; Text = "Introdução" },
    @{ Row = 33; Num = "5).4)"; Text = "Implementação" },
    @{ Row = 34; Num = "5).5)"; Text = "Conclusão" }
)
foreach ($item in $rows32to34) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = 0
    $ws.Range("B$r").Value = $item.Num
    $ws.Range("C$r").Value = $item.Text

    $ws.Range("B$r").Borders.LineStyle = 1

    $ws.Range("C$r`:D$r").Merge()
    $ws.Range("C$r`:D$r").Borders.LineStyle = 1
    $ws.Range("C$r`:D$r").HorizontalAlignment = -4131   # xlLeft
}

# ---------------------------------------------------------------------
# 4) Move the saved selection marker to I15 (cosmetic, matches the diff).
# ---------------------------------------------------------------------
$ws.Range("I15").Select()
